$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Header row: new columns I and J (set values after formatting so they are not overwritten)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-18
$data = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(8, 9)
    5  = @(4, 5)
    6  = @(8, 8)
    7  = @(9, 9)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(4, 5)
    16 = @(4, 4)
    17 = @(8, 8)
    18 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
